$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column P, which inherits column O's formatting (row by row)
# exactly as the author did when extending the table with a new year (2022).
$ws.Columns("P").Insert()

$ws.Range("P4").Value = 2022
$ws.Range("P5").Value = 96.969944810665083
$ws.Range("P6").Value = 96.173557859042035
$ws.Range("P7").Value = 62.289845326160055
$ws.Range("P8").Value = 100
$ws.Range("P9").Value = 100
$ws.Range("P10").Value = "-"
$ws.Range("P11").Value = 100
$ws.Range("P12").Value = 58.090784503861151
$ws.Range("P13").Value = 100
$ws.Range("P14").Value = 100

$ws.Range("Q4").Select()
